$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 holds the single comparable facility record. Replace it with the
# new facility's data (previously "Regina Health Center", now
# "Villa At Evergreen Park,The"), and update the distance / bed count /
# occupancy / star rating figures accordingly.

$ws.Range("B2").Value = "Villa At Evergreen Park,The"
$ws.Range("C2").Value = "10124 SOUTH KEDZIE"
$ws.Range("D2").Value = "Evergreen Park"
$ws.Range("E2").Value = "IL"
$ws.Range("F2").Value = "60805"
$ws.Range("G2").Value = "7089077000"
$ws.Range("I2").Value = 242
$ws.Range("J2").Value = "19.7 mi"
$ws.Range("L2").Value = 0.5826446280991735
$ws.Range("Q2").Value = "2"
